# Update the "想去人数" (want-to-go count) column F values for matching
# rows on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 19
$ws1.Range("F4").Value  = 1372
$ws1.Range("F5").Value  = 311
$ws1.Range("F7").Value  = 10665
$ws1.Range("F12").Value = 700
$ws1.Range("F13").Value = 12036
$ws1.Range("F14").Value = 12478
$ws1.Range("F16").Value = 119

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 19
$ws4.Range("F5").Value  = 1372
$ws4.Range("F6").Value  = 311
$ws4.Range("F8").Value  = 10665
$ws4.Range("F13").Value = 700
$ws4.Range("F14").Value = 12036
$ws4.Range("F15").Value = 12479
$ws4.Range("F17").Value = 119
